$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing the existing row 4 (and below) down.
$ws.Rows("4:4").Insert()

# Fill in the new row 4 with the new contact's phone number.
$ws.Range("A4").Value = 990033942
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = ""
